{"js": "// Lab2 report edits:\n//  1. Merge the \"One Integer operation ... One floating operation ...\" paragraph\n//     into a single clean run (drops the gramStart/gramEnd proofErr markers\n//     that wrapped the mid-sentence \"One\").\n//  2. Update the optimisation paragraph: \"one warp (32 threads)\" -> \"half of a\n//     warp (16 threads)\" and append \"without multiple threads trying to read\n//     from the same data block\" before the final period.\n//  3. Merge the \"I didn't have to manually allocate ...\" paragraph into a\n//     single clean run (drops the spellStart/spellEnd proofErr markers that\n//     wrapped \"device_vector\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  if (para.text.indexOf(\"One Integer operation to find the index of the vector\") !== -1) {\n    para.clear();\n    para.insertText(\n      \"One Integer operation to find the index of the vector, One floating operation to add the values at a given index in the vectors\",\n      Word.InsertLocation.start\n    );\n  } else if (para.text.indexOf(\"I didn\\u2019t have to manually allocate and free global memory\") !== -1) {\n    para.clear();\n    para.insertText(\n      \"I didn\\u2019t have to manually allocate and free global memory, this was handled by the device_vector allocation. I didn\\u2019t have to write configure the grid and block dimensions. I didn\\u2019t have to write a vector addition kernel.\",\n      Word.InsertLocation.start\n    );\n  }\n}\nawait context.sync();\n\n// Replace \"one warp (32 threads)\" with \"half of a warp (16 threads)\" \u2014 keep the\n// search match one character short of the run boundary (i.e. drop the trailing\n// space) so the shorter replacement text doesn't bleed into the next run.\nconst warpResults = body.search(\"no more than one warp (32 threads)\", { matchCase: false });\nwarpResults.load(\"text\");\nawait context.sync();\nfor (const result of warpResults.items) {\n  result.insertText(\n    \"no more than half of a warp (16 threads)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// Extend \"on separate multiprocessors.\" with the extra clause about avoiding\n// multiple threads reading from the same data block.\nconst multiprocResults = body.search(\"on separate multiprocessors.\", { matchCase: false });\nmultiprocResults.load(\"text\");\nawait context.sync();\nfor (const result of multiprocResults.items) {\n  result.insertText(\n    \"on separate multiprocessors without multiple threads trying to read from the same data block.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Lab2 report edits (\"added correct optimisation to lab2 report\"):\n#  1. Merge the \"One Integer operation ... One floating operation ...\" text\n#     into a single clean run (drops the gramStart/gramEnd proofing marks\n#     that wrapped the mid-sentence \"One\").\n#  2. Update the optimisation paragraph: \"one warp (32 threads)\" -> \"half of\n#     a warp (16 threads)\" and append \"without multiple threads trying to\n#     read from the same data block\" before the final period.\n#  3. Merge the \"I didn't have to manually allocate ...\" text into a single\n#     clean run (drops the spellStart/spellEnd proofing marks that wrapped\n#     \"device_vector\").\n#\n# wdReplaceAll = 2 (used positionally as the last Find.Execute argument so\n# every occurrence in scope is replaced).\n\n$d = $word.ActiveDocument\n$apostrophe = [char]0x2019\n\n# 1. Re-type the \"One Integer operation...\" sentence as a single run.\n$d.Content.Find.Execute(\n    \"One Integer operation to find the index of the vector, One floating operation to add the values at a given index in the vectors\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"One Integer operation to find the index of the vector, One floating operation to add the values at a given index in the vectors\",\n    2\n) | Out-Null\n\n# 2. Warp size: one warp (32 threads) -> half of a warp (16 threads).\n$d.Content.Find.Execute(\n    \"no more than one warp (32 threads)\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"no more than half of a warp (16 threads)\",\n    2\n) | Out-Null\n\n# 3. Add the clause about avoiding multiple threads reading the same data block.\n$d.Content.Find.Execute(\n    \"on separate multiprocessors.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"on separate multiprocessors without multiple threads trying to read from the same data block.\",\n    2\n) | Out-Null\n\n# 4. Re-type the \"I didn't have to manually allocate...\" paragraph as a single run.\n$thrustText = \"I didn${apostrophe}t have to manually allocate and free global memory, this was handled by the device_vector allocation. I didn${apostrophe}t have to write configure the grid and block dimensions. I didn${apostrophe}t have to write a vector addition kernel.\"\n$d.Content.Find.Execute(\n    $thrustText,\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    $thrustText,\n    2\n) | Out-Null\n"}
